$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 181.83333
$ws.Range("I12").Value = 174
$ws.Range("K12").Value = 174
$ws.Range("M12").Value = -4
$ws.Range("H15").Value = 321.63635
$ws.Range("I15").Value = 321.63635
$ws.Range("K15").Value = 964.90905
$ws.Range("M15").Value = -795.90905
$ws.Range("H33").Value = 322.66666
$ws.Range("I33").Value = 282.94116
$ws.Range("J33").Value = 998
$ws.Range("K33").Value = 282.94116
$ws.Range("L33").Value = 998
$ws.Range("M33").Value = -53.94116000000002
$ws.Range("N33").Value = -1456
$ws.Range("H39").Value = 74.85714
$ws.Range("J39").Value = 300
$ws.Range("L39").Value = 900
$ws.Range("N39").Value = -1492
$ws.Range("H51").Value = 151999.83
$ws.Range("J51").Value = 102399.8
$ws.Range("L51").Value = 102399.8
$ws.Range("N51").Value = -103367.8
$ws.Range("H103").Value = 722.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 722.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2167.5
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -3339.5
$ws.Range("H138").Value = 2555.4546
$ws.Range("I138").Value = 1282
$ws.Range("J138").Value = 3033
$ws.Range("K138").Value = 3846
$ws.Range("L138").Value = 9099
$ws.Range("M138").Value = 1294
$ws.Range("N138").Value = -19379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 498.5
$ws.Range("I4").Value = 498
$ws.Range("K4").Value = 498
$ws.Range("M4").Value = -382
$ws.Range("H32").Value = 2299.5
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -713
$ws.Range("H96").Value = 50172
$ws.Range("J96").Value = 50172
$ws.Range("L96").Value = 50172
$ws.Range("N96").Value = -55664
$ws.Range("H97").Value = 1273.0588
$ws.Range("I97").Value = 860.1429000000001
$ws.Range("J97").Value = 3200
$ws.Range("K97").Value = 860.1429000000001
$ws.Range("L97").Value = 3200
$ws.Range("M97").Value = -364.1429000000001
$ws.Range("N97").Value = -4192
$ws.Range("H124").Value = 80999.5
$ws.Range("I124").Value = 79999
$ws.Range("J124").Value = 82000
$ws.Range("K124").Value = 79999
$ws.Range("L124").Value = 82000
$ws.Range("M124").Value = -75089
$ws.Range("N124").Value = -91820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 32000
$ws.Range("I88").Value = 13000
$ws.Range("J88").Value = 41500
$ws.Range("K88").Value = 13000
$ws.Range("L88").Value = 41500
$ws.Range("M88").Value = -12594
$ws.Range("N88").Value = -42312
$ws.Range("H91").Value = 32000
$ws.Range("I91").Value = 13000
$ws.Range("J91").Value = 41500
$ws.Range("K91").Value = 13000
$ws.Range("L91").Value = 41500
$ws.Range("M91").Value = -11596
$ws.Range("N91").Value = -44308
$ws.Range("H94").Value = 2175.85
$ws.Range("I94").Value = 1792.3334
$ws.Range("K94").Value = 1792.3334
$ws.Range("M94").Value = -1341.3334
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 10000
$ws.Range("K99").Value = 10000
$ws.Range("M99").Value = -8502
$ws.Range("H107").Value = 4437.4375
$ws.Range("I107").Value = 4437.4375
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4437.4375
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2517.4375
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 150000
$ws.Range("J112").Value = 150000
$ws.Range("L112").Value = 150000
$ws.Range("N112").Value = -152954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 80000
$ws.Range("J43").Value = 80000
$ws.Range("L43").Value = 80000
$ws.Range("N43").Value = -80368
$ws.Range("H88").Value = 23999
$ws.Range("J88").Value = 23999
$ws.Range("L88").Value = 23999
$ws.Range("N88").Value = -24811
$ws.Range("H91").Value = 23999
$ws.Range("J91").Value = 23999
$ws.Range("L91").Value = 23999
$ws.Range("N91").Value = -26807
$ws.Range("H101").Value = 80000
$ws.Range("J101").Value = 80000
$ws.Range("L101").Value = 80000
$ws.Range("N101").Value = -86490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.375
$ws.Range("J2").Value = 43.4
$ws.Range("L2").Value = 260.4
$ws.Range("N2").Value = -486.4
$ws.Range("H4").Value = 2759.9714
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 12000
$ws.Range("N4").Value = -12224
$ws.Range("H34").Value = 627.2222
$ws.Range("I34").Value = 356.42856
$ws.Range("J34").Value = 1575
$ws.Range("K34").Value = 1069.28568
$ws.Range("L34").Value = 4725
$ws.Range("M34").Value = -985.28568
$ws.Range("N34").Value = -4893
$ws.Range("H38").Value = 532.94116
$ws.Range("I38").Value = 100
$ws.Range("J38").Value = 560
$ws.Range("K38").Value = 300
$ws.Range("L38").Value = 1680
$ws.Range("M38").Value = 47
$ws.Range("N38").Value = -2374
$ws.Range("H55").Value = 2038.3846
$ws.Range("J55").Value = 4066.6667
$ws.Range("L55").Value = 12200.0001
$ws.Range("N55").Value = -12554.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3593.2
$ws.Range("I97").Value = 4238.75
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 4238.75
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -3742.75
$ws.Range("N97").Value = -2003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2501500
$ws.Range("I14").Value = 2501500
$ws.Range("K14").Value = 2501500
$ws.Range("M14").Value = -2501328
$ws.Range("H22").Value = 3999.75
$ws.Range("I22").Value = 3000.5
$ws.Range("K22").Value = 3000.5
$ws.Range("M22").Value = -2705.5
$ws.Range("H27").Value = 3999.75
$ws.Range("I27").Value = 3000.5
$ws.Range("K27").Value = 3000.5
$ws.Range("M27").Value = -2893.5
$ws.Range("H55").Value = 535.1429000000001
$ws.Range("I55").Value = 383
$ws.Range("J55").Value = 649.25
$ws.Range("K55").Value = 383
$ws.Range("L55").Value = 649.25
$ws.Range("M55").Value = -210
$ws.Range("N55").Value = -995.25
$ws.Range("H68").Value = 1901.2
$ws.Range("I68").Value = 1834.6666
$ws.Range("J68").Value = 2001
$ws.Range("K68").Value = 1834.6666
$ws.Range("L68").Value = 2001
$ws.Range("M68").Value = -1085.6666
$ws.Range("N68").Value = -3499
$ws.Range("H71").Value = 1901.2
$ws.Range("I71").Value = 1834.6666
$ws.Range("J71").Value = 2001
$ws.Range("K71").Value = 9173.333000000001
$ws.Range("L71").Value = 10005
$ws.Range("M71").Value = -5429.333000000001
$ws.Range("N71").Value = -17493
$ws.Range("H93").Value = 3166.2856
$ws.Range("I93").Value = 3166.2856
$ws.Range("K93").Value = 3166.2856
$ws.Range("M93").Value = -1918.2856
$ws.Range("H100").Value = 2166.6667
$ws.Range("I100").Value = 2166.6667
$ws.Range("K100").Value = 2166.6667
$ws.Range("M100").Value = -1625.6667
$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 2900
$ws.Range("K132").Value = 8700
$ws.Range("M132").Value = -6170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 40000
$ws.Range("J29").Value = 40000
$ws.Range("L29").Value = 40000
$ws.Range("N29").Value = -40580
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H81").Value = 2773.1667
$ws.Range("I81").Value = 2637.8
$ws.Range("J81").Value = 3450
$ws.Range("K81").Value = 5275.6
$ws.Range("L81").Value = 6900
$ws.Range("M81").Value = -4214.6
$ws.Range("N81").Value = -9022
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H84").Value = 2773.1667
$ws.Range("I84").Value = 2637.8
$ws.Range("J84").Value = 3450
$ws.Range("K84").Value = 26378
$ws.Range("L84").Value = 34500
$ws.Range("M84").Value = -21074
$ws.Range("N84").Value = -45108
$ws.Range("H96").Value = 1050
$ws.Range("I96").Value = 900
$ws.Range("K96").Value = 900
$ws.Range("M96").Value = 473
$ws.Range("H99").Value = 1432
$ws.Range("I99").Value = 1432
$ws.Range("K99").Value = 1432
$ws.Range("M99").Value = 1563
$ws.Range("H109").Value = 69999
$ws.Range("J109").Value = 69999
$ws.Range("L109").Value = 69999
$ws.Range("N109").Value = -72773
$ws.Range("H126").Value = 772.5714
$ws.Range("I126").Value = 772.5714
$ws.Range("K126").Value = 2317.7142
$ws.Range("M126").Value = 152.2857999999997
$ws.Range("H136").Value = 949.8333
$ws.Range("I136").Value = 949.8333
$ws.Range("K136").Value = 2849.4999
$ws.Range("M136").Value = -299.4998999999998

Write-Host "Applied all Golem_Profits market data updates."